$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-43) holds a "Förändrad" (last-changed) date serial number.
# Update every value from 45771 to 45772 (one day later),
# leaving formatting and all other cells untouched.
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45771) {
        $cell.Value = 45772
    }
}
